# Re-style the three tables (slides 14, 15 and 16) from the custom
# "Table_0" style ({53030788-BAD9-457A-B2C8-4625D192CFF3}) to the
# built-in "Light Style 1 - Accent 3" table style
# ({4FAC6652-6527-4AF7-BC96-5FBC2FFBE85C}).

$p = $ppt.ActivePresentation

$newStyleId = "{4FAC6652-6527-4AF7-BC96-5FBC2FFBE85C}"
$targetSlides = 14, 15, 16

foreach ($idx in $targetSlides) {
    $slide = $p.Slides.Item($idx)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shape = $slide.Shapes.Item($i)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle($newStyleId)
        }
    }
}
